# Generate Report for Handback
# Updates the zh-cn and de-de localization-status sheets with the
# handback result for 56294caa-3c2a-4fc6-b815-525f5fe865a2 (row 5):
#  - records the Latest Target File / Latest Handback File / Latest Handback DateTime
#  - flags that the handed-back file version is not the latest and records the error detail
#  - widens the "Error Detail" column so the message is readable

$wb = $excel.ActiveWorkbook

$latestMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b7c0a296a317a4bbfcfd4373cd4994eeb927811e/e2e/56294caa-3c2a-4fc6-b815-525f5fe865a2.md"
$mdDisplay = "56294caa-3c2a-4fc6-b815-525f5fe865a2.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7f5fd63aeb11d79c1c42f71262801f13d0794818/e2e/56294caa-3c2a-4fc6-b815-525f5fe865a2.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b7c0a296a317a4bbfcfd4373cd4994eeb927811e/e2e/56294caa-3c2a-4fc6-b815-525f5fe865a2.md."
# Same blue (FF6495ED) already used by the other hyperlink cells in these sheets (OLE BGR order)
$hyperlinkColor = 15570276

# ----- zh-cn sheet -----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Columns.Item(16).ColumnWidth = 39.17

$wsZh.Range("I5").Value = $mdDisplay
$wsZh.Hyperlinks.Add($wsZh.Range("I5"), $latestMdUrl, [Type]::Missing, [Type]::Missing, $mdDisplay)
$wsZh.Range("I5").Font.Color = $hyperlinkColor
$wsZh.Range("J5").Value = "56294caa-3c2a-4fc6-b815-525f5fe865a2.b2cd98ced85cfd9cdf8443389ca692d323894499.zh-cn.xlf"
$wsZh.Range("K5").Value = "2016-10-20 00:29:58"
$wsZh.Range("P5").Value = $errorDetail

# ----- de-de sheet -----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Columns.Item(16).ColumnWidth = 39.17

$wsDe.Range("I5").Value = $mdDisplay
$wsDe.Hyperlinks.Add($wsDe.Range("I5"), $latestMdUrl, [Type]::Missing, [Type]::Missing, $mdDisplay)
$wsDe.Range("I5").Font.Color = $hyperlinkColor
$wsDe.Range("J5").Value = "56294caa-3c2a-4fc6-b815-525f5fe865a2.b2cd98ced85cfd9cdf8443389ca692d323894499.de-de.xlf"
$wsDe.Range("K5").Value = "2016-10-20 00:30:17"
$wsDe.Range("P5").Value = $errorDetail

Write-Output "Report generated for handback of 56294caa-3c2a-4fc6-b815-525f5fe865a2"
